$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3 (pushes the existing rows 3-8 down to 4-9)
$ws.Rows.Item(3).Insert()

# Populate the newly inserted row 3 with the new weekly price record
$ws.Range("A3").Value = 4
$ws.Range("B3").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C3").Value = "Los Lagos"
$ws.Range("D3").Value = 44525
$ws.Range("E3").Value = 10
$ws.Range("F3").Value = 100112012
$ws.Range("G3").Value = "Espinaca"
$ws.Range("H3").Value = "Sin especificar"
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 20
$ws.Range("K3").Value = 9000
$ws.Range("L3").Value = 9000
$ws.Range("M3").Value = 9000
$ws.Range("N3").Value = "$/cuna 10 kilos"
$ws.Range("O3").Value = "Región Metropolitana"
$ws.Range("P3").Value = 900
$ws.Range("Q3").Value = 10
$ws.Range("R3").Value = "Hortaliza"
